$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 32258064
$ws.Range("I86").Value = 32258064
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 32258064
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -32256941
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 32258064
$ws.Range("I89").Value = 32258064
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 161290320
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -161284704
$ws.Range("N89").ClearContents()

$ws.Range("H135").Value = 6684
$ws.Range("I135").Value = 1590
$ws.Range("J135").Value = 11268.6
$ws.Range("K135").Value = 14310
$ws.Range("L135").Value = 101417.4
$ws.Range("M135").Value = -11775
$ws.Range("N135").Value = -106487.4

$ws.Range("H140").Value = 71510.95
$ws.Range("I140").Value = 73333.336
$ws.Range("J140").Value = 69081.11
$ws.Range("K140").Value = 73333.336
$ws.Range("L140").Value = 69081.11
$ws.Range("M140").Value = -68153.336
$ws.Range("N140").Value = -79441.11

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4137.619
$ws.Range("I32").Value = 2307.46
$ws.Range("J32").Value = 11176.692
$ws.Range("K32").Value = 2307.46
$ws.Range("L32").Value = 11176.692
$ws.Range("M32").Value = -2020.46
$ws.Range("N32").Value = -11750.692

$ws.Range("H61").Value = 11914852
$ws.Range("I61").Value = 15633025
$ws.Range("J61").Value = 16695.5
$ws.Range("K61").Value = 15633025
$ws.Range("L61").Value = 16695.5
$ws.Range("M61").Value = -15632813
$ws.Range("N61").Value = -17119.5

$ws.Range("H74").Value = 29413280
$ws.Range("I74").Value = 45455988
$ws.Range("J74").Value = 1646.5
$ws.Range("K74").Value = 45455988
$ws.Range("L74").Value = 1646.5
$ws.Range("M74").Value = -45455114
$ws.Range("N74").Value = -3394.5

$ws.Range("H77").Value = 29413280
$ws.Range("I77").Value = 45455988
$ws.Range("J77").Value = 1646.5
$ws.Range("K77").Value = 227279940
$ws.Range("L77").Value = 8232.5
$ws.Range("M77").Value = -227275572
$ws.Range("N77").Value = -16968.5

$ws.Range("H136").Value = 11914852
$ws.Range("I136").Value = 15633025
$ws.Range("J136").Value = 16695.5
$ws.Range("K136").Value = 46899075
$ws.Range("L136").Value = 50086.5
$ws.Range("M136").Value = -46896525
$ws.Range("N136").Value = -55186.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1470.8572
$ws.Range("I134").Value = 1382.8334
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 4148.5002
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -1613.5002
$ws.Range("N134").Value = -11067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3061.9565
$ws.Range("I31").Value = 2877.0715
$ws.Range("J31").Value = 5003.25
$ws.Range("K31").Value = 2877.0715
$ws.Range("L31").Value = 5003.25
$ws.Range("M31").Value = -2582.0715
$ws.Range("N31").Value = -5593.25

$ws.Range("H34").Value = 3061.9565
$ws.Range("I34").Value = 2877.0715
$ws.Range("J34").Value = 5003.25
$ws.Range("K34").Value = 2877.0715
$ws.Range("L34").Value = 5003.25
$ws.Range("M34").Value = -2675.0715
$ws.Range("N34").Value = -5407.25

$ws.Range("H80").Value = 98999.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 98999.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 98999.5
$ws.Range("N80").Value = -101245.5

$ws.Range("H83").Value = 98999.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 98999.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 296998.5
$ws.Range("N83").Value = -308230.5

$ws.Range("H132").Value = 2296.75
$ws.Range("I132").Value = 1552.2
$ws.Range("J132").Value = 3537.6667
$ws.Range("K132").Value = 4656.6
$ws.Range("L132").Value = 10613.0001
$ws.Range("M132").Value = -2126.6
$ws.Range("N132").Value = -15673.0001

$ws.Range("H134").Value = 3572.7334
$ws.Range("I134").Value = 3054.0908
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 9162.2724
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -6627.2724
$ws.Range("N134").Value = -20067

$ws.Range("H141").Value = 85590.64999999999
$ws.Range("I141").Value = 41999.5
$ws.Range("J141").Value = 91402.8
$ws.Range("K141").Value = 41999.5
$ws.Range("L141").Value = 91402.8
$ws.Range("M141").Value = -36819.5
$ws.Range("N141").Value = -101762.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 403.66666
$ws.Range("I22").Value = 294.44446
$ws.Range("J22").Value = 731.3333
$ws.Range("K22").Value = 883.33338
$ws.Range("L22").Value = 2193.9999
$ws.Range("M22").Value = -714.33338
$ws.Range("N22").Value = -2531.9999

$ws.Range("H27").Value = 403.66666
$ws.Range("I27").Value = 294.44446
$ws.Range("J27").Value = 731.3333
$ws.Range("K27").Value = 883.33338
$ws.Range("L27").Value = 2193.9999
$ws.Range("M27").Value = -781.33338
$ws.Range("N27").Value = -2397.9999

$ws.Range("H33").Value = 385.42856
$ws.Range("I33").Value = 283.16666
$ws.Range("J33").Value = 999
$ws.Range("K33").Value = 1698.99996
$ws.Range("L33").Value = 5994
$ws.Range("M33").Value = -1415.99996
$ws.Range("N33").Value = -6560

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H38").Value = 185.22728
$ws.Range("I38").Value = 56.76923
$ws.Range("J38").Value = 370.77777
$ws.Range("K38").Value = 170.30769
$ws.Range("L38").Value = 1112.33331
$ws.Range("M38").Value = 176.69231
$ws.Range("N38").Value = -1806.33331

$ws.Range("H40").Value = 2894.1428
$ws.Range("I40").Value = 41.333332
$ws.Range("J40").Value = 5033.75
$ws.Range("K40").Value = 165.333328
$ws.Range("L40").Value = 20135
$ws.Range("M40").Value = -96.33332799999999
$ws.Range("N40").Value = -20273

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H42").Value = 9998.4
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 9998.4
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 29995.2
$ws.Range("N42").Value = -31063.2

$ws.Range("H43").Value = 12700
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12700
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 38100
$ws.Range("N43").Value = -38328

$ws.Range("H44").Value = 524.9
$ws.Range("I44").Value = 393.75
$ws.Range("J44").Value = 1049.5
$ws.Range("K44").Value = 1181.25
$ws.Range("L44").Value = 3148.5
$ws.Range("M44").Value = -783.25
$ws.Range("N44").Value = -3944.5

$ws.Range("H46").Value = 3891.8572
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 5049
$ws.Range("K46").Value = 2997
$ws.Range("L46").Value = 15147
$ws.Range("M46").Value = -2906
$ws.Range("N46").Value = -15329

$ws.Range("H47").Value = 8253.286
$ws.Range("I47").Value = 4800
$ws.Range("J47").Value = 8828.833000000001
$ws.Range("K47").Value = 14400
$ws.Range("L47").Value = 26486.499
$ws.Range("M47").Value = -13969
$ws.Range("N47").Value = -27348.499

$ws.Range("H48").Value = 7875
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 7875
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 23625
$ws.Range("N48").Value = -24125

$ws.Range("H49").Value = 6206.6665
$ws.Range("I49").Value = 303
$ws.Range("J49").Value = 7387.4
$ws.Range("K49").Value = 909
$ws.Range("L49").Value = 22162.2
$ws.Range("M49").Value = -753
$ws.Range("N49").Value = -22474.2

$ws.Range("H56").Value = 7132.3
$ws.Range("I56").Value = 7132.3
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7132.3
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -6602.3

$ws.Range("H62").Value = 18994
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 21659.666
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 64978.99800000001
$ws.Range("M62").Value = -8314
$ws.Range("N62").Value = -66350.99800000001

$ws.Range("H65").Value = 18994
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 21659.666
$ws.Range("K65").Value = 27000
$ws.Range("L65").Value = 194936.994
$ws.Range("M65").Value = -23568
$ws.Range("N65").Value = -201800.994

$ws.Range("H74").Value = 27122.875
$ws.Range("I74").Value = 26996
$ws.Range("J74").Value = 27141
$ws.Range("K74").Value = 80988
$ws.Range("L74").Value = 81423
$ws.Range("M74").Value = -79927
$ws.Range("N74").Value = -83545

$ws.Range("H77").Value = 27122.875
$ws.Range("I77").Value = 26996
$ws.Range("J77").Value = 27141
$ws.Range("K77").Value = 242964
$ws.Range("L77").Value = 244269
$ws.Range("M77").Value = -237660
$ws.Range("N77").Value = -254877

$ws.Range("H108").Value = 14961.417
$ws.Range("I108").Value = 4941.3335
$ws.Range("J108").Value = 24981.5
$ws.Range("K108").Value = 14824.0005
$ws.Range("L108").Value = 74944.5
$ws.Range("M108").Value = -11944.0005
$ws.Range("N108").Value = -80704.5

$ws.Range("H109").Value = 11457.8
$ws.Range("I109").Value = 2432.6667
$ws.Range("J109").Value = 24995.5
$ws.Range("K109").Value = 7298.000100000001
$ws.Range("L109").Value = 74986.5
$ws.Range("M109").Value = -6258.000100000001
$ws.Range("N109").Value = -77066.5

$ws.Range("H120").Value = 23225
$ws.Range("I120").Value = 5562.5
$ws.Range("J120").Value = 35000
$ws.Range("K120").Value = 16687.5
$ws.Range("L120").Value = 105000
$ws.Range("M120").Value = -11849.5
$ws.Range("N120").Value = -114676

$ws.Range("H126").Value = 12500
$ws.Range("I126").Value = 12500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 37500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -32560

$ws.Range("H133").Value = 16777.8
$ws.Range("I133").Value = 20705.666
$ws.Range("J133").Value = 15965.138
$ws.Range("K133").Value = 62116.99800000001
$ws.Range("L133").Value = 47895.414
$ws.Range("M133").Value = -57056.99800000001
$ws.Range("N133").Value = -58015.414

$ws.Range("H139").Value = 1004892.9
$ws.Range("I139").Value = 1671004.9
$ws.Range("J139").Value = 5725
$ws.Range("K139").Value = 5013014.699999999
$ws.Range("L139").Value = 17175
$ws.Range("M139").Value = -5007874.699999999
$ws.Range("N139").Value = -27455

$ws.Range("H141").Value = 9577.727999999999
$ws.Range("I141").Value = 5106.875
$ws.Range("J141").Value = 21500
$ws.Range("K141").Value = 15320.625
$ws.Range("L141").Value = 64500
$ws.Range("M141").Value = -10140.625
$ws.Range("N141").Value = -74860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5779.4634
$ws.Range("I102").Value = 5277
$ws.Range("J102").Value = 7566
$ws.Range("K102").Value = 5277
$ws.Range("L102").Value = 7566
$ws.Range("M102").Value = -3655
$ws.Range("N102").Value = -10810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1251146
$ws.Range("I82").Value = 1645756.9
$ws.Range("J82").Value = 1544.8334
$ws.Range("K82").Value = 1645756.9
$ws.Range("L82").Value = 1544.8334
$ws.Range("M82").Value = -1645395.9
$ws.Range("N82").Value = -2266.8334

$ws.Range("H85").Value = 1251146
$ws.Range("I85").Value = 1645756.9
$ws.Range("J85").Value = 1544.8334
$ws.Range("K85").Value = 1645756.9
$ws.Range("L85").Value = 1544.8334
$ws.Range("M85").Value = -1644508.9
$ws.Range("N85").Value = -4040.8334

$ws.Range("H139").Value = 85000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 85000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6570
$ws.Range("I62").Value = 5457.143
$ws.Range("J62").Value = 9166.666999999999
$ws.Range("K62").Value = 5457.143
$ws.Range("L62").Value = 9166.666999999999
$ws.Range("M62").Value = -4833.143
$ws.Range("N62").Value = -10414.667

$ws.Range("H65").Value = 6570
$ws.Range("I65").Value = 5457.143
$ws.Range("J65").Value = 9166.666999999999
$ws.Range("K65").Value = 27285.715
$ws.Range("L65").Value = 45833.335
$ws.Range("M65").Value = -24165.715
$ws.Range("N65").Value = -52073.335

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
